$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Facility" column (column C) from the report.
$ws.Range("C1").EntireColumn.Delete()

$ws.Range("C1").Select()
